$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Column H: "وضعیت پارکینگ" -> "پارکینگ سندی"
$ws.Cells.Item(1, 8).Value = "پارکینگ سندی"
# Column F: "مبلغ شارژ ( ریال)" -> "مبلغ شارژ"
$ws.Cells.Item(1, 6).Value = "مبلغ شارژ"

# --- Data rows: column F (charge amount) becomes a formatted text string ---
# Group 1: E = "دوخوابه" rows -> old numeric 4780000 becomes text " 4,780,000 ریال "
# These are rows 2-149 excluding 138-141 (which belong to the "چهار خوابه" group)
for ($r = 2; $r -le 149; $r++) {
    if ($r -ge 138 -and $r -le 141) { continue }
    $ws.Cells.Item($r, 6).NumberFormat = "@"
    $ws.Cells.Item($r, 6).Value = " 4,780,000 ریال "
}

# Group 2: E = "سه خوابه" rows (150-213) -> old numeric 5190000 becomes text "  ریال 5,190,000"
for ($r = 150; $r -le 213; $r++) {
    $ws.Cells.Item($r, 6).Value = "  ریال 5,190,000"
}

# Group 3: E = "چهار خوابه" rows (138-141) -> old numeric 5710000 becomes text "5,710,000 ریال"
for ($r = 138; $r -le 141; $r++) {
    $ws.Cells.Item($r, 6).Value = "5,710,000 ریال"
}

# Column G header: "بدهی شما تا پایان  مرداد ماه 1404 (ریال)" -> "بدهی تا پایان  مرداد ماه 1404"
$ws.Cells.Item(1, 7).Value = "بدهی تا پایان  مرداد ماه 1404"

# --- Selection moves from M8:M9 to G10 ---
[void]$ws.Range("G10").Select()
